# Update the InstallTracker worksheet:
#  - Correct the "Actual" install counts for 8/27 and 8/28 (rows 97-98)
#  - Fill in newly tracked "Actual" install counts for 8/29 - 9/1 (rows 99-102)
#  - Extend the dependent Daily (D) / Average (E) shared formulas to match
#  - Restore the view to the top of the sheet / a different active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrected / new "Actual" values in column C ---------------------------
$actual = @{
    97  = 1942
    98  = 1992
    99  = 2070
    100 = 2113
    101 = 2157
    102 = 2213
}

foreach ($row in 97..102) {
    $cell = $ws.Cells.Item($row, 3)   # column C
    $cell.Value = $actual[$row]
    $cell.NumberFormat = "#,##0"
}

# --- Extend the "Daily" (D) shared formula down through row 102 -----------
foreach ($row in 99..102) {
    $dcell = $ws.Cells.Item($row, 4)
    $dcell.Formula = "=C$row-C" + ($row - 1)
    $dcell.NumberFormat = "0"
}

# --- Extend the "Average" (E) shared formula down through row 102 ---------
foreach ($row in 99..102) {
    $ecell = $ws.Cells.Item($row, 5)
    $ecell.Formula = "=(C$row-C" + ($row - 7) + ")/7"
    $ecell.NumberFormat = "0"
}

# --- Restore the frozen-pane view / selection ------------------------------
$ws.Activate()
$ws.Range("K33").Select()
